# Fix "unknown" tag values in the xlm (per-language tag-stats sheets).
# The commit message states: "Fixed unknown in xlm (rest has not changed)".
# Three sheets (Hebrew, Algerian, Maltese) had the Tag column (A/F/K/P) for a
# handful of low-frequency tag rows (18-21) mis-assigned; this corrects them.

$wb = $excel.ActiveWorkbook

# ---- Hebrew sheet ----
$ws = $wb.Worksheets.Item("Hebrew")
$ws.Range("F18").Value = "SYM"
$ws.Range("K18").Value = "SYM"
$ws.Range("A19").Value = "PART"
$ws.Range("F19").Value = "O"
$ws.Range("K19").Value = "O"
$ws.Range("P19").Value = "SYM"
$ws.Range("A20").Value = "SYM"
$ws.Range("F20").Value = "PART"
$ws.Range("K20").Value = "PART"
$ws.Range("P20").Value = "O"
$ws.Range("A21").Value = "O"
$ws.Range("F21").Value = "X"
$ws.Range("P21").Value = "PART"

# ---- Algerian sheet ----
$ws = $wb.Worksheets.Item("Algerian")
$ws.Range("K18").Value = "SYM"
$ws.Range("F19").Value = "NUM"
$ws.Range("K19").Value = "O"
$ws.Range("A20").Value = "SYM"
$ws.Range("P20").Value = "SYM"
$ws.Range("A21").Value = "O"
$ws.Range("F21").Value = "O"
$ws.Range("K21").Value = "NUM"
$ws.Range("P21").Value = "O"

# ---- Maltese sheet ----
$ws = $wb.Worksheets.Item("Maltese")
$ws.Range("A20").Value = "O"
$ws.Range("F20").Value = "O"
$ws.Range("K20").Value = "O"
$ws.Range("P20").Value = "O"
$ws.Range("A21").Value = "_"
$ws.Range("F21").Value = "_"
$ws.Range("K21").Value = "_"
$ws.Range("P21").Value = "_"
